# Fruta / hortaliza, semanal
# Update Fecha (D), Calidad (I), Volumen (J), Precio minimo (K),
# Precio maximo (L), Precio promedio ponderado (M) and Precio $/Kg (P)
# for the weekly consolidated price records (rows 2-11 and 13-18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = 44474; I = "Segunda"; J = 200;  K = 600; L = 700; M = 650; P = 650 },
    @{ Row = 3;  D = 44229; I = "Segunda"; J = 760;  K = 550; L = 600; M = 575; P = 575 },
    @{ Row = 4;  D = 44210; I = "Segunda"; J = 900;  K = 600; L = 700; M = 650; P = 650 },
    @{ Row = 5;  D = 44174; I = "Segunda"; J = 800;  K = 450; L = 500; M = 475; P = 475 },
    @{ Row = 6;  D = 44174; I = "Tercera"; J = 1200; K = 250; L = 350; M = 300; P = 300 },
    @{ Row = 7;  D = 44573; I = "Tercera"; J = 800;  K = 600; L = 650; M = 625; P = 625 },
    @{ Row = 8;  D = 44658; I = "Segunda"; J = 1000; K = 600; L = 650; M = 625; P = 625 },
    @{ Row = 9;  D = 44253; I = "Segunda"; J = 1000; K = 800; L = 900; M = 850; P = 850 },
    @{ Row = 10; D = 44253; I = "Tercera"; J = 800;  K = 600; L = 700; M = 650; P = 650 },
    @{ Row = 11; D = 44201; I = "Segunda"; J = 500;  K = 800; L = 900; M = 850; P = 850 },
    @{ Row = 13; D = 44544; I = "Primera"; J = 1000; K = 600; L = 650; M = 625; P = 625 },
    @{ Row = 14; D = 44278; I = "Segunda"; J = 700;  K = 600; L = 700; M = 650; P = 650 },
    @{ Row = 15; D = 44278; I = "Tercera"; J = 400;  K = 500; L = 600; M = 550; P = 550 },
    @{ Row = 16; D = 44245; I = "Primera"; J = 800;  K = 850; L = 900; M = 875; P = 875 },
    @{ Row = 17; D = 44245; I = "Segunda"; J = 1000; K = 750; L = 800; M = 775; P = 775 },
    @{ Row = 18; D = 44267; I = "Tercera"; J = 400;  K = 500; L = 600; M = 550; P = 550 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("D$r").Value = $u.D
    $ws.Range("I$r").Value = $u.I
    $ws.Range("J$r").Value = $u.J
    $ws.Range("K$r").Value = $u.K
    $ws.Range("L$r").Value = $u.L
    $ws.Range("M$r").Value = $u.M
    $ws.Range("P$r").Value = $u.P
}
